$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price column (D) cells we touch stay text, not auto-converted to numbers/dates,
# matching the original workbook where Price values are stored as text.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "24.899.13"
$ws.Range("E2").Value = "  +1.86%  "
$ws.Range("D3").Value = "1.671.69"
$ws.Range("E3").Value = "  +1.00%  "
$ws.Range("D4").Value = "1.004"
$ws.Range("E4").Value = "  +0.22%  "
$ws.Range("D5").Value = "329.01"
$ws.Range("E5").Value = "  +6.82%  "
$ws.Range("D6").Value = "1.001"
$ws.Range("E6").Value = "  +0.20%  "
$ws.Range("D7").Value = "0.3650"
$ws.Range("E7").Value = "  +0.83%  "
$ws.Range("D8").Value = "46.49"
$ws.Range("E8").Value = "  -1.76%  "
$ws.Range("D9").Value = "0.3239"
$ws.Range("E9").Value = "  -1.24%  "
$ws.Range("D10").Value = "1.140"
$ws.Range("E10").Value = "  +1.10%  "
$ws.Range("D11").Value = "0.07048"
$ws.Range("E11").Value = "  +1.43%  "
$ws.Range("D12").Value = "1.002"
$ws.Range("E12").Value = "  +0.38%  "
$ws.Range("D13").Value = "6.061"
$ws.Range("E13").Value = "  +1.72%  "
$ws.Range("D14").Value = "19.54"
$ws.Range("E14").Value = "  +1.09%  "
$ws.Range("D15").Value = "1.674.58"
$ws.Range("E15").Value = "  +1.41%  "
$ws.Range("D16").Value = "6.607"
$ws.Range("E16").Value = "  -0.33%  "
$ws.Range("D17").Value = "0.00001044"
$ws.Range("E17").Value = "  +0.26%  "
$ws.Range("D18").Value = "0.06578"
$ws.Range("E18").Value = "  +0.88%  "
$ws.Range("D19").Value = "1.001"
$ws.Range("E19").Value = "  +0.29%  "
$ws.Range("D20").Value = "78.68"
$ws.Range("E20").Value = "  +2.77%  "
$ws.Range("D21").Value = "15.82"
$ws.Range("E21").Value = "  +0.59%  "
$ws.Range("D22").Value = "5.910"
$ws.Range("E22").Value = "  -0.20%  "
$ws.Range("D23").Value = "12.95"
$ws.Range("E23").Value = "  +2.86%  "
$ws.Range("D24").Value = "24.928.40"
$ws.Range("E24").Value = "  +2.16%  "
$ws.Range("D25").Value = "2.438"
$ws.Range("E25").Value = "  +0.38%  "
$ws.Range("D26").Value = "2.386"
$ws.Range("E26").Value = "  +1.64%  "
$ws.Range("D27").Value = "148.10"
$ws.Range("E27").Value = "  +1.03%  "
$ws.Range("D28").Value = "18.73"
$ws.Range("E28").Value = "  +2.01%  "
$ws.Range("D29").Value = "1.863.57"
$ws.Range("E29").Value = "  +1.26%  "
$ws.Range("E30").Value = "  +0.94%  "
$ws.Range("D31").Value = "1.180"
$ws.Range("E31").Value = "  -0.07%  "
$ws.Range("D32").Value = "4.067"
$ws.Range("E32").Value = "  +0.63%  "
$ws.Range("D33").Value = "5.750"
$ws.Range("E33").Value = "  +1.82%  "
$ws.Range("D34").Value = "0.08447"
$ws.Range("E34").Value = "  +1.35%  "
$ws.Range("D35").Value = "1.647"
$ws.Range("E35").Value = "  -1.51%  "
$ws.Range("D36").Value = "12.26"
$ws.Range("E36").Value = "  -0.35%  "
$ws.Range("D37").Value = "5.152"
$ws.Range("E37").Value = "  -1.43%  "
$ws.Range("B38").Value = "TrustWalletToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D38").Value = "1.234"
$ws.Range("E38").Value = "  +2.40%  "
$ws.Range("B39").Value = "VeChain"
$ws.Range("C39").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D39").Value = "0.02240"
$ws.Range("E39").Value = "  +1.53%  "
$ws.Range("B40").Value = "Algorand"
$ws.Range("C40").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D40").Value = "0.2088"
$ws.Range("E40").Value = "  +1.75%  "
$ws.Range("B41").Value = "Hedera"
$ws.Range("C41").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D41").Value = "0.05999"
$ws.Range("E41").Value = "  -0.83%  "
$ws.Range("D42").Value = "8.212"
$ws.Range("E42").Value = "  +0.22%  "
$ws.Range("D43").Value = "0.9992"
$ws.Range("E43").Value = "  +0.07%  "
$ws.Range("D44").Value = "0.5932"
$ws.Range("E44").Value = "  +1.53%  "
$ws.Range("D45").Value = "13.68"
$ws.Range("E45").Value = "  +8.11%  "
$ws.Range("E46").Value = "  +2.89%  "
$ws.Range("D47").Value = "0.5724"
$ws.Range("E47").Value = "  +2.66%  "
$ws.Range("D48").Value = "124.71"
$ws.Range("E48").Value = "  +2.16%  "
$ws.Range("D49").Value = "1.958"
$ws.Range("E49").Value = "  +0.79%  "
$ws.Range("D50").Value = "0.07006"
$ws.Range("E50").Value = "  +1.55%  "
$ws.Range("D51").Value = "1.185"
$ws.Range("E51").Value = "  +2.52%  "

# Restore default (unstyled) formatting on the Price column so cells match the
# workbook's original style (no explicit number format), while keeping the text values.
$ws.Range("D2:D51").Style = "Normal"
